$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194..233 down to 195..234.
$ws.Rows("194:194").Insert()

# Populate the newly inserted row 194 with its values.
# Columns A,B,C,E,F,G,H,I,J,N,O,Q,R keep the same values the old row 194 had
# (now located at row 195); only D,K,L,M,P change.
$ws.Range("A194").Value = 7
$ws.Range("B194").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C194").Value = "Ñuble"
$ws.Range("D194").Value = 44694
$ws.Range("E194").Value = 16
$ws.Range("F194").Value = 100112003
$ws.Range("G194").Value = "Ajo"
$ws.Range("H194").Value = "Chino"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 100
$ws.Range("K194").Value = 20000
$ws.Range("L194").Value = 21000
$ws.Range("M194").Value = 20500
$ws.Range("N194").Value = "`$/caja 10 kilos"
$ws.Range("O194").Value = "China"
$ws.Range("P194").Value = 2050
$ws.Range("Q194").Value = 10
$ws.Range("R194").Value = "Hortaliza"

# Make sure the D194 cell keeps the date number format used by the rest of column D.
$ws.Range("D194").NumberFormat = $ws.Range("D195").NumberFormat
